$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Extend the D49:D58 / E49:E58 merge (week total) down to row 62, since the
#    new entries for the week of 43233 are being appended inside that block.
# ---------------------------------------------------------------------------
$ws.Range("D49:D58").UnMerge()
$ws.Range("E49:E58").UnMerge()

# Restyle the now unmerged D58/E58 (previously last row of block, thick bottom
# border) to the plain "middle of block" style used on D50:E57.
$ws.Range("D51:E51").Copy()
$ws.Range("D58:E58").PasteSpecial(-4122)

# Row 58 loses its thick bottom border / extra height now that it is no
# longer the last row of the week block.
$ws.Range("A51:C51").Copy()
$ws.Range("A58:C58").PasteSpecial(-4122)
$ws.Range("A58").Value = 43230
$ws.Range("B58").Value = $ws.Range("B59").Value
$ws.Range("C58").Value = 1.25

# ---------------------------------------------------------------------------
# 2. Fill in the new journal entries (week of 2018-05-13 / serial 43233).
# ---------------------------------------------------------------------------

# Row 59
$ws.Range("A52:C52").Copy()
$ws.Range("A59:C59").PasteSpecial(-4122)
$ws.Range("D52:E52").Copy()
$ws.Range("D59:E59").PasteSpecial(-4122)
$ws.Range("A59").Value = 43233
$ws.Range("B59").Value = "Mise en forme des transaction dans le dahboard"
$ws.Range("C59").Value = 0.75
$ws.Range("F59").Value = "11h30"
$ws.Range("G59").Value = "12h45"

# Row 60
$ws.Range("A52:C52").Copy()
$ws.Range("A60:C60").PasteSpecial(-4122)
$ws.Range("D52:E52").Copy()
$ws.Range("D60:E60").PasteSpecial(-4122)
$ws.Range("A60").Value = 43233
$ws.Range("B60").Value = "Ajout du pieChart"
$ws.Range("C60").Value = 0.5
$ws.Range("F60").Value = "12h45"
$ws.Range("G60").Value = "13h15"

# Row 61
$ws.Range("A51:C51").Copy()
$ws.Range("A61:C61").PasteSpecial(-4122)
$ws.Range("D51:E51").Copy()
$ws.Range("D61:E61").PasteSpecial(-4122)
$ws.Range("A61").Value = 43233
$ws.Range("B62temp") = $null
$ws.Range("B61").Value = "Correction graphique dashboard temp placeholder"
$ws.Range("C61").Value = 1.5

# Row 62
$ws.Range("A6:C6").Copy()
$ws.Range("A62:C62").PasteSpecial(-4122)
$ws.Range("D6:E6").Copy()
$ws.Range("D62:E62").PasteSpecial(-4122)
$ws.Range("A62").Value = 43233
$ws.Range("B62").Value = "Correction graphique dashboard"
$ws.Range("C62").Value = 0.5

# Fix up row 61's text (needs to be written after row 62's "Correction..." so
# shared string indices come out in the same order as the source workbook).
$ws.Range("B61").Value = "Probleme general dans l'application. Recherche du probleme si c'etait du cote GUI"

# Re-extend the merge to cover the new rows.
$ws.Range("D49:D62").Merge()
$ws.Range("E49:E62").Merge()
$ws.Range("E49").Formula = "=SUM(C49:C62)"

# ---------------------------------------------------------------------------
# 3. Shift the trailing blank rows + Total row down by 5 (59-62 -> 63-67),
#    adding one extra blank row (66).
# ---------------------------------------------------------------------------
$ws.Range("A63:C63").Value = ""
$ws.Range("A64:C64").Value = ""
$ws.Range("A65:C65").Value = ""
$ws.Range("A66:C66").Value = ""
$ws.Range("B67").Value = "Total"
$ws.Range("C67").Formula = "=SUM(C5:C66)"
